$wb = $excel.ActiveWorkbook

# OFF sheet - offensive target depth stats, Home row (row 2)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 618
$wsOff.Range("C2").Value = 431
$wsOff.Range("D2").Value = 147
$wsOff.Range("E2").Value = 69
$wsOff.Range("F2").Value = 7

# DEF sheet - defensive target depth stats, Home row (row 2)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 577
$wsDef.Range("C2").Value = 400
$wsDef.Range("D2").Value = 131
$wsDef.Range("E2").Value = 46
$wsDef.Range("F2").Value = 13
